# Replace the FALSE() boolean-formula cells in I2:I27 with the literal text
# string "False" (the sample sheet used a custom "TRUE"/"FALSE" number format
# over a boolean formula; the edit turns those into plain text cells reading
# "False").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 27; $r++) {
    $c = $ws.Cells.Item($r, 9)    # column I
    # Format the cell as Text first so the literal value we paste in isn't
    # re-interpreted as a boolean.
    $c.NumberFormat = "@"
    # Put the text in via a formula, then flatten the formula down to its
    # literal value (copy / paste-special-values) so the cell ends up holding
    # plain text "False" rather than a live formula or an auto-coerced
    # boolean.
    $c.Formula = "=""False"""
    $c.Copy()
    $c.PasteSpecial(-4163)  # xlPasteValues
}

# Mirror the author's final selection (I2:I27, active cell I2).
$ws.Range("I2:I27").Select() | Out-Null
